# Change "dress" -> "blouse" throughout the experiment protocol document.
#
# Commit message: "changed dress to blouse to be less ambiguously
# relevant to women."
#
# This touches three spots:
#   1. the target-item line "9: dress"                 -> "10: blouse"
#   2. "...in order, "pants", "dress", "rings"..."      -> "...", "blouse", ..."
#   3. "DRESS: "Please select dress.""                  -> "BLOUSE: "Please select blouse.""
#
# Word's auto-managed "_GoBack" bookmark (which marks the location of the
# most recent edit) also moves from the end of the document to right
# after the word "blouse" typed in spot 3 - the last edit made.

$d = $word.ActiveDocument

function Get-ParaByText($needle) {
    $count = $d.Paragraphs.Count
    for ($i = 1; $i -le $count; $i++) {
        $p = $d.Paragraphs($i)
        if ($p.Range.Text -like $needle) {
            return $p
        }
    }
    return $null
}

$wdFindContinue = 1
$wdReplaceOne = 1

# --- 1. "9: dress" -> "10: blouse" -------------------------------------
$p1 = Get-ParaByText("9: dress*")
$r = $p1.Range
$r.Find.Execute("9", $true, $true, $false, $false, $false, $true, $wdFindContinue, $false, "10", $wdReplaceOne)
$r = $p1.Range
$r.Find.Execute("dress", $true, $true, $false, $false, $false, $true, $wdFindContinue, $false, "blouse", $wdReplaceOne)

# --- 2. "..., "pants", "dress", "rings", ..." -> "..., "blouse", ..." --
$p2 = Get-ParaByText("*order of targets*")
$r = $p2.Range
$r.Find.Execute("dress", $true, $true, $false, $false, $false, $true, $wdFindContinue, $false, "blouse", $wdReplaceOne)

# --- 3. "DRESS: "Please select dress."" -> "BLOUSE: "Please select blouse."" --
$p3 = Get-ParaByText("DRESS:*")
$r = $p3.Range
$r.Find.Execute("DRESS", $true, $true, $false, $false, $false, $true, $wdFindContinue, $false, "BLOUSE", $wdReplaceOne)

$p3 = Get-ParaByText("BLOUSE:*")
$r = $p3.Range
$r.Find.Execute("dress", $true, $true, $false, $false, $false, $true, $wdFindContinue, $false, "blouse", $wdReplaceOne)
# $r now covers exactly the just-inserted "blouse" (Find.Execute leaves the
# range spanning the replacement text, mirroring real Word COM behaviour).
$r.Collapse(0)  # wdCollapseEnd - collapse to the point right after "blouse"

# --- 4. Relocate the "_GoBack" bookmark to that point -------------------
$old = $d.Bookmarks("_GoBack")
$old.Delete()
$d.Bookmarks.Add("_GoBack", $r)

Write-Output "done"
